# Edit the "Requisition and Issue Slip" workbook to match the target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block: Division / Responsibility Center / RIS No. ---
$ws.Range("B8").Value = "KMD Regular"
$ws.Range("H8").Value = "KMD ISSP PuRD 200000100002"
$ws.Range("H9").Value = ""

# --- Item rows (12 & 13): Stock No. cleared, Unit/Description replaced, quantities updated, Remarks cleared ---
$ws.Range("A12").Value = ""
$ws.Range("B12").Value = "pcs"
$ws.Range("C12").Value = "as scdsfd safsdfdf dsaf"
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = "R"
$ws.Range("G12").Value = 10
$ws.Range("H12").Value = ""

$ws.Range("A13").Value = ""
$ws.Range("B13").Value = "pcs"
$ws.Range("C13").Value = "as scdsfd safsdfdf dsaf"
$ws.Range("D13").Value = 5
$ws.Range("E13").Value = "R"
$ws.Range("G13").Value = 10
$ws.Range("H13").Value = ""

# --- Purpose ---
$ws.Range("B15").Value = "sadasdasdasdasds"

# --- Signatory block: names / designations cleared, dates updated ---
$ws.Range("D22").Value = ""
$ws.Range("F22").Value = ""
$ws.Range("H22").Value = ""

$ws.Range("C23").Value = ""
$ws.Range("D23").Value = ""
$ws.Range("F23").Value = ""
$ws.Range("H23").Value = ""

$ws.Range("C24").Value = "12/19/2022"
$ws.Range("D24").Value = "12/19/2022"
$ws.Range("F24").Value = "12/19/2022"
$ws.Range("H24").Value = "12/19/2022"
